# Apply crypto price/volume update (GitHub Actions scheduled refresh)
# Updates D (Price) and E (Volume(1h)) columns for many rows, and
# for a few rows the underlying coin (B/C/D/E) changed position entirely.
# For D-column values that look like plain numbers, force text format first
# so Excel keeps them as text (matching the original inline-string cells)
# instead of silently converting them to floating point numbers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '59.418.65'
$ws.Range("E2").Value = '  +4.52%  '
$ws.Range("D3").Value = '2.600.07'
$ws.Range("E3").Value = '  +2.44%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '522.29'
$ws.Range("E5").Value = '  +1.72%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '141.62'
$ws.Range("E6").Value = '  +1.12%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.998'
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.568'
$ws.Range("E8").Value = '  +2.54%  '
$ws.Range("D9").Value = '2.624.75'
$ws.Range("E9").Value = '  +3.28%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '6.53'
$ws.Range("E10").Value = '  +0.52%  '
$ws.Range("E11").Value = '  +2.23%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.333'
$ws.Range("E12").Value = '  +3.20%  '
$ws.Range("E13").Value = '  +2.61%  '
$ws.Range("D14").Value = '3.061.96'
$ws.Range("E14").Value = '  +2.56%  '
$ws.Range("D15").Value = '59.359.54'
$ws.Range("E15").Value = '  +4.36%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '20.63'
$ws.Range("E16").Value = '  +3.00%  '
$ws.Range("D17").Value = '2.643.31'
$ws.Range("E17").Value = '  +3.18%  '
$ws.Range("E18").Value = '  +0.80%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '340.23'
$ws.Range("E19").Value = '  +2.37%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.35'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.24'
$ws.Range("E22").Value = '  +7.33%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '0.998'
$ws.Range("E23").Value = '  -0.22%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '66.55'
$ws.Range("E24").Value = '  +3.85%  '
$ws.Range("E25").Value = '  +1.73%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.407'
$ws.Range("E26").Value = '  +1.81%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '1.00'
$ws.Range("E27").Value = '  +0.23%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.15'
$ws.Range("E28").Value = '  +4.12%  '
$ws.Range("E29").Value = '  -0.07%  '
$ws.Range("D30").Value = '0.0₃0734'
$ws.Range("E30").Value = '  -1.83%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '5.98'
$ws.Range("E31").Value = '  -4.17%  '
$ws.Range("E32").Value = '  +1.29%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '18.85'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '148.93'
$ws.Range("E34").Value = '  +0.60%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '4.03'
$ws.Range("E35").Value = '  +1.35%  '
$ws.Range("E36").Value = '  +0.97%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '36.39'
$ws.Range("E37").Value = '  +2.49%  '
$ws.Range("B38").Value = 'Fetch.AI'
$ws.Range("C38").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.841'
$ws.Range("E38").Value = '  +2.47%  '
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.46'
$ws.Range("E39").Value = '  +2.84%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.833'
$ws.Range("E40").Value = '  -0.69%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.57'
$ws.Range("E41").Value = '  +2.96%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '278.44'
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.998'
$ws.Range("E43").Value = '  -0.14%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '10.75'
$ws.Range("E44").Value = '  +1.34%  '
$ws.Range("B45").Value = 'Mantle'
$ws.Range("C45").Value = 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.593'
$ws.Range("E45").Value = '  +2.80%  '
$ws.Range("B46").Value = 'Stellar'
$ws.Range("C46").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.0957'
$ws.Range("E46").Value = '  +0.58%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '0.0523'
$ws.Range("E47").Value = '  +0.95%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '18.69'
$ws.Range("E48").Value = '  +1.06%  '
$ws.Range("D49").Value = '1.993.63'
$ws.Range("E49").Value = '  +1.49%  '
$ws.Range("B50").Value = 'VeChain'
$ws.Range("C50").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0221'
$ws.Range("E50").Value = '  +0.57%  '
$ws.Range("B51").Value = 'RenderToken'
$ws.Range("C51").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '4.54'
$ws.Range("E51").Value = '  +0.45%  '